$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(3408, 45565, 19267.03, 26.87, 4.13, 1.08),
    @(3409, 45566, 19274.99, 26.88, 4.13, 1.08),
    @(3410, 45568, 18857.48, 26.3, 4.04, 1.1),
    @(3411, 45569, 18679.27, 26.05, 4, 1.11),
    @(3412, 45572, 18433.79, 25.71, 3.95, 1.13),
    @(3413, 45573, 18683.35, 26.05, 4.01, 1.11),
    @(3414, 45574, 18724.98, 26.11, 4.01, 1.11),
    @(3415, 45575, 18717.78, 26.09, 4.01, 1.11),
    @(3416, 45576, 18717.57, 26.09, 4.01, 1.11),
    @(3417, 45579, 18816.93, 26.23, 4.03, 1.11),
    @(3418, 45580, 18794.6, 26.19, 4.02, 1.11),
    @(3419, 45581, 18732.72, 26.1, 4.01, 1.11),
    @(3420, 45582, 18505.3, 25.79, 3.96, 1.12),
    @(3421, 45583, 18569.04, 25.88, 3.97, 1.13),
    @(3422, 45586, 18448.92, 25.65, 3.95, 1.12),
    @(3423, 45587, 18136.33, 25.13, 3.88, 1.15),
    @(3424, 45588, 18139.15, 25.13, 3.88, 1.15),
    @(3425, 45589, 18102.38, 25.02, 3.87, 1.15),
    @(3426, 45590, 17894.11, 24.71, 3.83, 1.12),
    @(3427, 45593, 18010.12, 24.89, 3.85, 1.11),
    @(3428, 45594, 18133.32, 25.32, 3.88, 1.15),
    @(3429, 45595, 18062.75, 25.22, 3.86, 1.16),
    @(3430, 45596, 17972.34, 25.03, 3.84, 1.18),
    @(3431, 45597, 18067.56, 25.16, 3.86, 1.17),
    @(3432, 45600, 17834.72, 24.82, 3.81, 1.17),
    @(3433, 45601, 17973.7, 25, 3.84, 1.18),
    @(3434, 45602, 18247.49, 25.38, 3.9, 1.17),
    @(3435, 45603, 18060.24, 25.1, 3.86, 1.19),
    @(3436, 45604, 17956.13, 24.77, 3.84, 1.2),
    @(3437, 45607, 17935.66, 24.71, 3.83, 1.18),
    @(3438, 45608, 17726.21, 24.25, 3.79, 1.2),
    @(3439, 45609, 17428.02, 23.88, 3.73, 1.22),
    @(3440, 45610, 17455.28, 23.91, 3.73, 1.23),
    @(3441, 45614, 17416.49, 23.84, 3.72, 1.22),
    @(3442, 45615, 17488.16, 23.92, 3.74, 1.23),
    @(3443, 45617, 17336.24, 23.72, 3.7, 1.21),
    @(3444, 45618, 17680.96, 24.19, 3.78, 1.18),
    @(3445, 45621, 17931.05, 24.53, 3.83, 1.17),
    @(3446, 45622, 17925.48, 24.52, 3.83, 1.17),
    @(3447, 45623, 18023.56, 24.65, 3.84, 1.17),
    @(3448, 45624, 17858.39, 24.43, 3.8, 1.18),
    @(3449, 45625, 17993.43, 24.61, 3.83, 1.17),
    @(3450, 45628, 18112.16, 24.77, 3.85, 1.16),
    @(3451, 45629, 18255.95, 24.97, 3.88, 1.15),
    @(3452, 45630, 18309.27, 25.04, 3.89, 1.15),
    @(3453, 45631, 18464.03, 25.25, 3.93, 1.14),
    @(3454, 45632, 18481.73, 25.28, 3.93, 1.14),
    @(3455, 45635, 18455.81, 25.24, 3.92, 1.14),
    @(3456, 45636, 18466.69, 25.26, 3.93, 1.14),
    @(3457, 45637, 18495.95, 25.3, 3.93, 1.14),
    @(3458, 45638, 18428.45, 25.2, 3.92, 1.13),
    @(3459, 45639, 18533.35, 25.35, 3.94, 1.12),
    @(3460, 45642, 18514.65, 25.32, 3.94, 1.12),
    @(3461, 45643, 18298.4, 25.03, 3.89, 1.14),
    @(3462, 45644, 18168.13, 24.85, 3.86, 1.14),
    @(3463, 45645, 18008.16, 24.63, 3.83, 1.15),
    @(3464, 45646, 17659.24, 24.15, 3.75, 1.18),
    @(3465, 45649, 17759.39, 24.29, 3.78, 1.17),
    @(3466, 45650, 17742.85, 24.26, 3.77, 1.18),
    @(3467, 45652, 17769.69, 24.3, 3.78, 1.18),
    @(3468, 45653, 17774.22, 24.31, 3.78, 1.17),
    @(3469, 45656, 17707.01, 24.22, 3.76, 1.17),
    @(3470, 45657, 17704.45, 24.35, 3.78, 1.17),
    @(3471, 45658, 17777.76, 24.45, 3.8, 1.16),
    @(3472, 45659, 18071.69, 24.85, 3.86, 1.14),
    @(3473, 45660, 17973.24, 24.72, 3.84, 1.15),
    @(3474, 45663, 17605.21, 24.21, 3.76, 1.17),
    @(3475, 45664, 17676.19, 24.31, 3.78, 1.17),
    @(3476, 45665, 17612.77, 24.22, 3.76, 1.17),
    @(3477, 45666, 17468.83, 24.02, 3.73, 1.18),
    @(3478, 45667, 17306.21, 23.8, 3.7, 1.19),
    @(3479, 45670, 16902.47, 23.23, 3.61, 1.22),
    @(3480, 45671, 17081.58, 23.48, 3.65, 1.21)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}

Write-Host "done"
